$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "69.909.46"
$ws.Range("E2").Value = "  -0.46%  "
$ws.Range("D3").Value = "3.467.72"
$ws.Range("E3").Value = "  -1.58%  "
$ws.Range("E4").Value = "  -0.22%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "611.80"
$ws.Range("E5").Value = "  +1.14%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "168.15"
$ws.Range("E6").Value = "  -2.47%  "
$ws.Range("D7").Value = "3.463.62"
$ws.Range("E7").Value = "  -1.52%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.596"
$ws.Range("E8").Value = "  -2.30%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.999"
$ws.Range("E9").Value = "  -0.13%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.193"
$ws.Range("E10").Value = "  +0.40%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "7.14"
$ws.Range("E11").Value = "  -1.80%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.567"
$ws.Range("E12").Value = "  -2.71%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "44.55"
$ws.Range("E13").Value = "  -3.46%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.0000270"
$ws.Range("E14").Value = "  -1.75%  "
$ws.Range("D15").Value = "4.019.04"
$ws.Range("E15").Value = "  -1.81%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "8.22"
$ws.Range("E16").Value = "  -0.61%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "588.63"
$ws.Range("E17").Value = "  -2.83%  "
$ws.Range("B18").Value = "WrappedBTC"
$ws.Range("C18").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D18").Value = "69.960.59"
$ws.Range("E18").Value = "  -0.44%  "
$ws.Range("B19").Value = "WrappedEther"
$ws.Range("C19").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D19").Value = "3.463.31"
$ws.Range("E19").Value = "  -1.76%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.121"
$ws.Range("E20").Value = "  +1.09%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "17.22"
$ws.Range("E21").Value = "  +0.16%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.856"
$ws.Range("E22").Value = "  -1.73%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "8.98"
$ws.Range("E23").Value = "  -3.15%  "
$ws.Range("B24").Value = "Litecoin"
$ws.Range("C24").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "95.76"
$ws.Range("E24").Value = "  +0.01%  "
$ws.Range("B25").Value = "InternetComputer(DFINITY)"
$ws.Range("C25").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "15.22"
$ws.Range("E25").Value = "  -2.40%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.63"
$ws.Range("E26").Value = "  -2.39%  "
$ws.Range("E27").Value = "  +0.01%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.47"
$ws.Range("E28").Value = "  -4.26%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "33.01"
$ws.Range("E29").Value = "  -3.29%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "8.68"
$ws.Range("E30").Value = "  -3.36%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.89"
$ws.Range("E31").Value = "  -2.78%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.84"
$ws.Range("E32").Value = "  -5.77%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.25"
$ws.Range("E33").Value = "  -3.02%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.60"
$ws.Range("E34").Value = "  -4.79%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "583.27"
$ws.Range("E35").Value = "  -19.26%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0484"
$ws.Range("E36").Value = "  +1.89%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "10.64"
$ws.Range("E37").Value = "  -0.43%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0964"
$ws.Range("E38").Value = "  -3.40%  "
$ws.Range("E39").Value = "  +0.28%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "56.15"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.141"
$ws.Range("E41").Value = "  -1.05%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.20"
$ws.Range("E42").Value = "  -10.34%  "
$ws.Range("D43").Value = "3.252.31"
$ws.Range("E43").Value = "  -2.89%  "
$ws.Range("D44").Value = "0.0₃0702"
$ws.Range("E44").Value = "  +1.59%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.297"
$ws.Range("E45").Value = "  -5.57%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "30.94"
$ws.Range("E46").Value = "  -4.40%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.78"
$ws.Range("E47").Value = "  -4.41%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.41"
$ws.Range("E48").Value = "  -5.49%  "
$ws.Range("E49").Value = "  -2.73%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "133.75"
$ws.Range("E50").Value = "  +1.01%  "
